$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 previously only had a (blank/styled) A52 cell. Give it the same
# look-and-feel as the row above (B51) before filling in the new values,
# so the newly introduced B52 cell reuses the existing style rather than
# creating a brand new one.
$ws.Cells.Item(51, 2).Copy()
$ws.Cells.Item(52, 2).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(52, 1).Value = "test"
$ws.Cells.Item(52, 2).Value = "123"

# Move the active selection from A54 to A53
$ws.Range("A53").Select()
